# "add drop item list record"
#
# NPC.xlsx / Record_PosList sheet gains a new drop-item-list record: columns
# L:O used to hold the struct's field *type* ("float") on row 1 and the field
# *names* (X, Y, Z, StayTime) on row 2. This edit swaps them so row 1 holds
# the field names and row 2 holds the field types, and moves the matching
# header comments down to row 2 to match. View/selection state is also
# updated to reflect the Record_PosList sheet being the active tab.

$wb = $excel.ActiveWorkbook
$wsProperty = $wb.Worksheets.Item("Property")
$wsPosList  = $wb.Worksheets.Item("Record_PosList")

# --- Row 1 becomes the field names, row 2 becomes the field types ("float")
#     for the new X / Y / Z / StayTime drop position-list columns.
$wsPosList.Range("L1").Value = "X"
$wsPosList.Range("M1").Value = "Y"
$wsPosList.Range("N1").Value = "Z"
$wsPosList.Range("O1").Value = "StayTime"

$wsPosList.Range("L2").Value = "float"
$wsPosList.Range("M2").Value = "float"
$wsPosList.Range("N2").Value = "float"
$wsPosList.Range("O2").Value = "float"

# --- Re-create the comments describing columns L:O one row down (row 2
#     instead of row 1) to match the new header/type layout.
$wsPosList.Range("L1").Comment.Delete()
$wsPosList.Range("M1").Comment.Delete()
$wsPosList.Range("N1").Comment.Delete()
$wsPosList.Range("O1").Comment.Delete()

[void]$wsPosList.Range("L2").AddComment("强化等级")
[void]$wsPosList.Range("M2").AddComment("强化等级")
[void]$wsPosList.Range("N2").AddComment("强化等级")
[void]$wsPosList.Range("O2").AddComment("镶嵌宝石，逗号分隔")

# --- Selection / active-tab bookkeeping to match the saved view state:
#     Property's selection moves to J33 (and loses focus), Record_PosList
#     becomes the active tab with selection on O8.
[void]$wsProperty.Range("J33").Select()
[void]$wsPosList.Range("O8").Select()
